$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- P2 / P3: numeric "Ma so thue" values (style matches existing numeric-text col, no quote prefix) ---
$ws.Range("P2").Value = 8077026742
$ws.Range("J4").Copy()
$ws.Range("P2").PasteSpecial(-4122)

$ws.Range("P3").Value = 8002259299
$ws.Range("J4").Copy()
$ws.Range("P3").PasteSpecial(-4122)

# --- P4 .. P13: text "Ma so thue" values (kept as text, matching Q column's quote-prefixed text style) ---
# Two-pass paste-format trick: first pass flips the cell to the Text number format
# (so the literal digit string is entered/kept as text rather than re-interpreted as
# a number), second pass re-applies the same source's format so the final style index
# exactly matches the already-existing quote-prefixed style (no new style is created).

$ws.Range("Q4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = "8019373151"
$ws.Range("Q4").Copy()
$ws.Range("P4").PasteSpecial(-4122)

$ws.Range("Q4").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value = "8101391805"
$ws.Range("Q5").Copy()
$ws.Range("P5").PasteSpecial(-4122)

$ws.Range("Q6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = "0103417878"
$ws.Range("Q6").Copy()
$ws.Range("P6").PasteSpecial(-4122)

$ws.Range("Q7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P7").Value = "0103413841"
$ws.Range("Q7").Copy()
$ws.Range("P7").PasteSpecial(-4122)

$ws.Range("Q8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value = "8327096826"
$ws.Range("Q8").Copy()
$ws.Range("P8").PasteSpecial(-4122)

$ws.Range("Q9").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P9").Value = "8026501206"
$ws.Range("Q9").Copy()
$ws.Range("P9").PasteSpecial(-4122)

$ws.Range("Q10").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("P10").Value = "8009385976"
$ws.Range("Q10").Copy()
$ws.Range("P10").PasteSpecial(-4122)

$ws.Range("Q11").Copy()
$ws.Range("P11").PasteSpecial(-4122)
$ws.Range("P11").Value = "8012846722"
$ws.Range("Q11").Copy()
$ws.Range("P11").PasteSpecial(-4122)

$ws.Range("Q12").Copy()
$ws.Range("P12").PasteSpecial(-4122)
$ws.Range("P12").Value = "8009473608"
$ws.Range("Q12").Copy()
$ws.Range("P12").PasteSpecial(-4122)

$ws.Range("Q13").Copy()
$ws.Range("P13").PasteSpecial(-4122)
$ws.Range("P13").Value = "8539460058"
$ws.Range("Q13").Copy()
$ws.Range("P13").PasteSpecial(-4122)

# --- column width adjustments ---
$ws.Columns("O").ColumnWidth = 60.26953125
$ws.Columns("P").ColumnWidth = 35.6328125
$ws.Columns("R").ColumnWidth = 81.7265625

# --- view / selection state ---
$ws.Application.ActiveWindow.ScrollColumn = 16
$ws.Range("S7").Select()
